$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "29.336.89"
Set-TextValue $ws.Range("E2") "  +1.43%  "
Set-TextValue $ws.Range("D3") "1.903.90"
Set-TextValue $ws.Range("E3") "  +1.09%  "
Set-TextValue $ws.Range("E4") "  +0.01%  "
Set-TextValue $ws.Range("D5") "322.85"
Set-TextValue $ws.Range("E5") "  -2.01%  "
Set-TextValue $ws.Range("E6") "  +0.01%  "
Set-TextValue $ws.Range("D7") "0.4724"
Set-TextValue $ws.Range("E7") "  +3.04%  "
Set-TextValue $ws.Range("D8") "0.4031"
Set-TextValue $ws.Range("E8") "  -1.51%  "
Set-TextValue $ws.Range("B9") "Dogecoin"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D9") "0.08012"
Set-TextValue $ws.Range("E9") "  +0.65%  "
Set-TextValue $ws.Range("B10") "Polygon"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D10") "0.9917"
Set-TextValue $ws.Range("E10") "  -0.04%  "
Set-TextValue $ws.Range("B11") "Solana"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D11") "22.58"
Set-TextValue $ws.Range("E11") "  +4.50%  "
Set-TextValue $ws.Range("B12") "WrappedEther"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D12") "1.895.51"
Set-TextValue $ws.Range("E12") "  -0.22%  "
Set-TextValue $ws.Range("B13") "Polkadot"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D13") "5.864"
Set-TextValue $ws.Range("E13") "  -0.67%  "
Set-TextValue $ws.Range("B14") "Chainlink"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D14") "7.063"
Set-TextValue $ws.Range("E14") "  +0.09%  "
Set-TextValue $ws.Range("B15") "Litecoin"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D15") "89.26"
Set-TextValue $ws.Range("E15") "  +1.03%  "
Set-TextValue $ws.Range("B16") "BinanceUSD"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D16") "1.002"
Set-TextValue $ws.Range("E16") "  +0.02%  "
Set-TextValue $ws.Range("B17") "TRON"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D17") "0.06612"
Set-TextValue $ws.Range("E17") "  +0.46%  "
Set-TextValue $ws.Range("B18") "ShibaInu"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D18") "0.00001027"
Set-TextValue $ws.Range("E18") "  +0.19%  "
Set-TextValue $ws.Range("B19") "Avalanche"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D19") "17.54"
Set-TextValue $ws.Range("E19") "  +0.94%  "
Set-TextValue $ws.Range("B20") "Dai"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D20") "1.001"
Set-TextValue $ws.Range("E20") "  -0.06%  "
Set-TextValue $ws.Range("B21") "WrappedBTC"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D21") "29.354.36"
Set-TextValue $ws.Range("E21") "  +1.57%  "
Set-TextValue $ws.Range("B22") "Uniswap"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D22") "5.513"
Set-TextValue $ws.Range("E22") "  +2.03%  "
Set-TextValue $ws.Range("B23") "Cosmos"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D23") "11.45"
Set-TextValue $ws.Range("E23") "  -0.13%  "
Set-TextValue $ws.Range("B24") "Toncoin"
Set-TextValue $ws.Range("C24") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D24") "2.201"
Set-TextValue $ws.Range("E24") "  +0.21%  "
Set-TextValue $ws.Range("B25") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C25") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D25") "2.118.93"
Set-TextValue $ws.Range("E25") "  -0.01%  "
Set-TextValue $ws.Range("B26") "Monero"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D26") "154.25"
Set-TextValue $ws.Range("E26") "  -1.25%  "
Set-TextValue $ws.Range("B27") "EthereumClassic"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D27") "19.72"
Set-TextValue $ws.Range("E27") "  +0.95%  "
Set-TextValue $ws.Range("B28") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D28") "6.036"
Set-TextValue $ws.Range("E28") "  +10.48%  "
Set-TextValue $ws.Range("B29") "LidoDAOToken"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D29") "2.089"
Set-TextValue $ws.Range("E29") "  +0.58%  "
Set-TextValue $ws.Range("B30") "BitcoinCash"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D30") "117.56"
Set-TextValue $ws.Range("E30") "  +0.15%  "
Set-TextValue $ws.Range("B31") "ImmutableX"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D31") "1.065"
Set-TextValue $ws.Range("E31") "  +4.37%  "
Set-TextValue $ws.Range("B32") "Stellar"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D32") "0.09444"
Set-TextValue $ws.Range("E32") "  +1.37%  "
Set-TextValue $ws.Range("B33") "ARBITRUM"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D33") "1.411"
Set-TextValue $ws.Range("E33") "  +0.99%  "
Set-TextValue $ws.Range("B34") "HuobiToken"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D34") "3.553"
Set-TextValue $ws.Range("E34") "  +0.95%  "
Set-TextValue $ws.Range("B35") "Filecoin"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D35") "5.353"
Set-TextValue $ws.Range("E35") "  +1.34%  "
Set-TextValue $ws.Range("B36") "Hedera"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D36") "0.06083"
Set-TextValue $ws.Range("E36") "  +0.59%  "
Set-TextValue $ws.Range("B37") "VeChain"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D37") "0.02244"
Set-TextValue $ws.Range("E37") "  +0.88%  "
Set-TextValue $ws.Range("B38") "TrustWalletToken"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D38") "1.173"
Set-TextValue $ws.Range("E38") "  +0.12%  "
Set-TextValue $ws.Range("B39") "FraxShare"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D39") "8.050"
Set-TextValue $ws.Range("E39") "  -3.07%  "
Set-TextValue $ws.Range("B40") "TheSandbox"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D40") "0.5814"
Set-TextValue $ws.Range("E40") "  +0.68%  "
Set-TextValue $ws.Range("B41") "RenderToken"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D41") "2.494"
Set-TextValue $ws.Range("E41") "  +11.40%  "
Set-TextValue $ws.Range("B42") "Algorand"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D42") "0.1829"
Set-TextValue $ws.Range("E42") "  +0.42%  "
Set-TextValue $ws.Range("B43") "Aptos"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D43") "10.07"
Set-TextValue $ws.Range("E43") "  +0.31%  "
Set-TextValue $ws.Range("B44") "WEMIXToken"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D44") "1.271"
Set-TextValue $ws.Range("E44") "  +1.04%  "
Set-TextValue $ws.Range("B45") "Cronos"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D45") "0.07692"
Set-TextValue $ws.Range("E45") "  +2.37%  "
Set-TextValue $ws.Range("B46") "EnergySwap"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D46") "12.18"
Set-TextValue $ws.Range("E46") "  +1.98%  "
Set-TextValue $ws.Range("B47") "Decentraland"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D47") "0.5483"
Set-TextValue $ws.Range("E47") "  +0.72%  "
Set-TextValue $ws.Range("B48") "NEARProtocol"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D48") "1.903"
Set-TextValue $ws.Range("E48") "  +0.33%  "
Set-TextValue $ws.Range("B49") "Quant"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D49") "113.28"
Set-TextValue $ws.Range("E49") "  +1.77%  "
Set-TextValue $ws.Range("B50") "WOONetwork"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-TextValue $ws.Range("D50") "0.2931"
Set-TextValue $ws.Range("E50") "  +5.27%  "
Set-TextValue $ws.Range("D51") "43.64"
Set-TextValue $ws.Range("E51") "  -3.40%  "
